$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $value) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $value
    $c.Style = "Normal"
}

# New data row (row 2)
Set-TextValue "A2" "Dustin"
Set-TextValue "B2" "12"
Set-TextValue "C2" "Male"
Set-TextValue "D2" "123445278"
Set-TextValue "E2" "2398578957"
Set-TextValue "F2" "899900581424"
Set-TextValue "G2" "[]"

# New column G header: numeric 0, matching the bordered header style used by A1:F1
$ws.Range("F1").Copy()
$ws.Range("G1").PasteSpecial(-4122)
$ws.Range("G1").Value = 0

# Page margins (Excel Normal margins)
$ws.PageSetup.LeftMargin = 0.7 * 72
$ws.PageSetup.RightMargin = 0.7 * 72
$ws.PageSetup.TopMargin = 0.75 * 72
$ws.PageSetup.BottomMargin = 0.75 * 72
$ws.PageSetup.HeaderMargin = 0.3 * 72
$ws.PageSetup.FooterMargin = 0.3 * 72

# Drop the stale cell selection left over from the previous session
$ws.Range("A1").Select() | Out-Null
